# Append run: 2025-10-05 06:24 JST
# The scraper re-ran, dropped the stale/no-longer-relevant listings that were
# previously in rows 2,3,5,7,8 (and their now-unused "skill concept" tags in
# column H), and kept only the two still-open listings (previously rows 4 and
# 6), refreshing their "fetched at" timestamp and moving them up to rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop every existing data row (2-8) in one shot; header row 1 stays put.
$ws.Range("A2:A8").EntireRow.Delete()

# The hyperlink objects don't follow the row shift, so clear them out -
# they'll be recreated below for the two rows that survive.
$ws.Hyperlinks.Delete()

# Row 2: 【急募】FXトレード履歴を基にしたEA作成依頼 (was row 4)
$ws.Range("A2").Value = "2025-10-05 06:24:47"
$ws.Range("B2").Value = "【急募】FXトレード履歴を基にしたEA作成依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5406904"
$ws.Range("G2").Value = 18

# Row 3: 【SalesIQ活用】CRMと連携したリード獲得方法を教えてください (was row 6)
$ws.Range("A3").Value = "2025-10-05 06:24:47"
$ws.Range("B3").Value = "【SalesIQ活用】CRMと連携したリード獲得方法を教えてください"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "~ 5,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5400402"
$ws.Range("G3").Value = 10

# Re-create the URL hyperlinks + apply the built-in "Hyperlink" style,
# matching the F-column formatting the other rows already had.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5406904")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5400402")
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# Column B (タイトル) got narrower in this revision: 47 -> 36 characters.
# ColumnWidth goes through Excel's char-width/pixel rounding, so 35.2 is the
# input that lands exactly on a stored width of 36.
$ws.Columns.Item(2).ColumnWidth = 35.2
